{"js": "// Apply \"hybrid bold + color\" highlighting to quantitative metrics\n// (percentages, dollar amounts, large numbers) in specific resume bullet\n// paragraphs, matching the author's commit: split each target run into\n// plain-text runs plus bold + color (2C3E50) runs around each metric.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Each entry: the *original* full paragraph text (used to uniquely locate\n// the paragraph) and the ordered list of metric substrings inside it that\n// should become bold + colored.\nconst targets = [\n  {\n    text:\n      \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    text:\n      \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00b14.2% to \\u00b12.1%\",\n    metrics: [\"87%\", \"71%\", \"\\u00b14.2%\", \"\\u00b12.1%\"],\n  },\n  {\n    text: \"\\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    metrics: [\"1,200\"],\n  },\n  {\n    text:\n      \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    metrics: [\"$400M\", \"$1B\"],\n  },\n  {\n    text: \"\\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    metrics: [\"87%\", \"71%\"],\n  },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map each desired full-text to the (first) matching paragraph object.\n// Paragraphs are processed in document order; each full-text is unique\n// among the body's paragraphs (confirmed against the source document), and\n// every paragraph is only used once, so indices naturally line up even\n// though one text is a strict prefix of another.\nconst used = new Set();\nfunction findParagraph(fullText) {\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (used.has(i)) continue;\n    if (paragraphs.items[i].text === fullText) {\n      used.add(i);\n      return paragraphs.items[i];\n    }\n  }\n  return null;\n}\n\nfor (const target of targets) {\n  const paragraph = findParagraph(target.text);\n  if (!paragraph) {\n    throw new Error(\"Could not locate paragraph: \" + target.text);\n  }\n\n  // Gather all the metric ranges first (search against the paragraph's\n  // original, still-unmodified text) before applying any formatting, so\n  // earlier in-place edits can't shift/invalidate later search hits.\n  const ranges = target.metrics.map((metric) =>\n    paragraph.search(metric, { matchCase: true, matchWholeWord: false })\n  );\n  ranges.forEach((r) => r.load(\"items\"));\n  await context.sync();\n\n  for (let i = 0; i < ranges.length; i++) {\n    const hits = ranges[i].items;\n    if (!hits || hits.length === 0) {\n      throw new Error(\"Metric not found: \" + target.metrics[i] + \" in \" + target.text);\n    }\n    const hit = hits[0];\n    hit.font.bold = true;\n    hit.font.color = HIGHLIGHT_COLOR;\n  }\n  await context.sync();\n}\n", "ps1": "# Apply \"hybrid bold + color\" highlighting to quantitative metrics\n# (percentages, dollar amounts, large numbers) in specific resume bullet\n# paragraphs, matching the author's commit: split each target run into\n# plain-text runs plus bold + color (2C3E50) runs around each metric.\n#\n# Word's Font.Color is a BGR-packed integer (0xBBGGRR), so RGB 2C3E50\n# becomes 0x503E2C = 5258796.\n\n$d = $word.ActiveDocument\n$highlightColor = 5258796\n\n$targets = @(\n    @{\n        Text    = \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Metrics = @(\"23%\", \"64%\")\n    },\n    @{\n        Text    = \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\"\n        Metrics = @(\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\")\n    },\n    @{\n        Text    = \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        Metrics = @(\"1,200\")\n    },\n    @{\n        Text    = \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        Metrics = @(\"`$400M\", \"`$1B\")\n    },\n    @{\n        Text    = \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        Metrics = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        Text    = \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        Metrics = @(\"87%\", \"71%\")\n    }\n)\n\n$paragraphCount = $d.Paragraphs.Count\n$usedIndices = @()\n\nfunction Find-ParagraphIndex($fullText, $usedIndices) {\n    for ($i = 1; $i -le $paragraphCount; $i++) {\n        if ($usedIndices -contains $i) {\n            continue\n        }\n        $p = $d.Paragraphs.Item($i)\n        $raw = $p.Range.Text\n        # Strip the trailing paragraph-mark character before comparing.\n        $candidate = $raw.Substring(0, $raw.Length - 1)\n        if ($candidate -eq $fullText) {\n            return $i\n        }\n    }\n    return -1\n}\n\nforeach ($target in $targets) {\n    $idx = Find-ParagraphIndex $target.Text $usedIndices\n    if ($idx -eq -1) {\n        throw \"Could not locate paragraph: \" + $target.Text\n    }\n    $usedIndices += $idx\n\n    foreach ($metric in $target.Metrics) {\n        $searchRange = $d.Paragraphs.Item($idx).Range\n        $hit = $searchRange.Find.Execute($metric)\n        if (-not $hit) {\n            throw \"Metric not found: \" + $metric + \" in paragraph \" + $idx\n        }\n        $searchRange.Font.Bold = 1\n        $searchRange.Font.Color = $highlightColor\n    }\n}\n\nWrite-Output \"done\"\n"}
